$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'
$ws.Range('B8').Value = 'Pabellón De Arteaga'
$ws.Range('B9').Value = 'Rincón De Romos'
$ws.Range('B10').Value = 'San José De Gracia'
$ws.Range('B32').Value = 'Comitán De Domínguez'
$ws.Range('B44').Value = 'Mazapa De Madero'
$ws.Range('B46').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B53').Value = 'San Cristóbal De Las Casas'
$ws.Range('B77').Value = 'Coyame Del Sotol'
$ws.Range('B85').Value = 'Guadalupe Y Calvo'
$ws.Range('B87').Value = 'Hidalgo Del Parral'
$ws.Range('B107').Value = 'San Francisco Del Oro'
$ws.Range('B113').Value = 'Valle De Zaragoza'
$ws.Range('B133').Value = 'San Juan De Sabinas'
$ws.Range('B142').Value = 'Villa De Álvarez'
$ws.Range('A144').Value = 'Ciudad De México'
$ws.Range('B148').Value = 'Cuajimalpa De Morelos'
$ws.Range('B163').Value = 'Coneto De Comonfort'
$ws.Range('B176').Value = 'Nombre De Dios'
$ws.Range('B179').Value = 'Pánuco De Coronado'
$ws.Range('B186').Value = 'San Juan De Guadalupe'
$ws.Range('B187').Value = 'San Juan Del Río'
$ws.Range('B188').Value = 'San Luis Del Cordero'
$ws.Range('B189').Value = 'San Pedro Del Gallo'
$ws.Range('A199').Value = 'Estado De México'
$ws.Range('B199').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B202').Value = 'Almoloya De Alquisiras'
$ws.Range('B203').Value = 'Almoloya Del Río'
$ws.Range('B206').Value = 'Atizapán De Zaragoza'
$ws.Range('B214').Value = 'Ecatepec De Morelos'
$ws.Range('B217').Value = 'Ixtapan De La Sal'
$ws.Range('B222').Value = 'Naucalpan De Juárez'
$ws.Range('B226').Value = 'San Felipe Del Progreso'
$ws.Range('B227').Value = 'San Martín De Las Pirámides'
$ws.Range('B233').Value = 'Tenango Del Aire'
$ws.Range('B236').Value = 'Tlalnepantla De Baz'
$ws.Range('B241').Value = 'Valle De Bravo'
$ws.Range('B242').Value = 'Villa De Allende'
$ws.Range('B252').Value = 'San Miguel De Allende'
$ws.Range('B253').Value = 'Apaseo El Alto'
$ws.Range('B254').Value = 'Apaseo El Grande'
$ws.Range('B259').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B263').Value = 'Jaral Del Progreso'
$ws.Range('B270').Value = 'Purísima Del Rincón'
$ws.Range('B273').Value = 'San Diego De La Unión'
$ws.Range('B275').Value = 'San Francisco Del Rincón'
$ws.Range('B277').Value = 'San Luis De La Paz'
$ws.Range('B279').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B281').Value = 'Silao De La Victoria'
$ws.Range('B286').Value = 'Valle De Santiago'
$ws.Range('B292').Value = 'Acapulco De Juárez'
$ws.Range('B296').Value = 'Atenango Del Río'
$ws.Range('B298').Value = 'Atoyac De Álvarez'
$ws.Range('B299').Value = 'Ayutla De Los Libres'
$ws.Range('B301').Value = 'Buenavista De Cuéllar'
$ws.Range('B302').Value = 'Chilapa De Álvarez'
$ws.Range('B303').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B305').Value = 'Coyuca De Benítez'
$ws.Range('B306').Value = 'Coyuca De Catalán'
$ws.Range('B309').Value = 'Cutzamala De Pinzón'
$ws.Range('B315').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B316').Value = 'Iguala De La Independencia'
$ws.Range('B318').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B320').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B323').Value = 'Mártir De Cuilapan'
$ws.Range('B330').Value = 'Taxco De Alarcón'
$ws.Range('B332').Value = 'Técpan De Galeana'
$ws.Range('B334').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B336').Value = 'Tixtla De Guerrero'
$ws.Range('B339').Value = 'Tlapa De Comonfort'
$ws.Range('B352').Value = 'Huejutla De Reyes'
$ws.Range('B355').Value = 'Jacala De Ledezma'
$ws.Range('B357').Value = 'Mineral Del Monte'
$ws.Range('B358').Value = 'Mixquiahuala De Juárez'
$ws.Range('B360').Value = 'Pachuca De Soto'
$ws.Range('B362').Value = 'Progreso De Obregón'
$ws.Range('B367').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B369').Value = 'Tezontepec De Aldama'
$ws.Range('B373').Value = 'Tula De Allende'
$ws.Range('B374').Value = 'Tulancingo De Bravo'
$ws.Range('B377').Value = 'Zacualtipán De Ángeles'
$ws.Range('B382').Value = 'Ahualulco De Mercado'
$ws.Range('B386').Value = 'Atotonilco El Alto'
$ws.Range('B387').Value = 'Autlán De Navarro'
$ws.Range('B391').Value = 'Cañadas De Obregón'
$ws.Range('B394').Value = 'Concepción De Buenos Aires'
$ws.Range('B399').Value = 'Encarnación De Díaz'
$ws.Range('B404').Value = 'Huejuquilla El Alto'
$ws.Range('B405').Value = 'Ixtlahuacán Del Río'
$ws.Range('B411').Value = 'Lagos De Moreno'
$ws.Range('B416').Value = 'Ojuelos De Jalisco'
$ws.Range('B420').Value = 'San Diego De Alejandría'
$ws.Range('B421').Value = 'San Juan De Los Lagos'
$ws.Range('B422').Value = 'San Juanito De Escobedo'
$ws.Range('B424').Value = 'San Martín De Bolaños'
$ws.Range('B426').Value = 'Santa María De Los Ángeles'
$ws.Range('B429').Value = 'Tamazula De Gordiano'
$ws.Range('B431').Value = 'Tepatitlán De Morelos'
$ws.Range('B432').Value = 'Tizapán El Alto'
$ws.Range('B438').Value = 'Unión De San Antonio'
$ws.Range('B443').Value = 'Yahualica De González Gallo'
$ws.Range('B444').Value = 'Zacoalco De Torres'
$ws.Range('B446').Value = 'Zapotlán Del Rey'
$ws.Range('B447').Value = 'Zapotlán El Grande'
$ws.Range('B465').Value = 'Cojumatlán De Régules'
$ws.Range('B512').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B531').Value = 'Coatlán Del Río'
$ws.Range('B539').Value = 'Puente De Ixtla'
$ws.Range('B542').Value = 'Tetela Del Volcán'
$ws.Range('B550').Value = 'Ixtlán Del Río'
$ws.Range('B564').Value = 'Ciénega De Flores'
$ws.Range('B572').Value = 'Mier Y Noriega'
$ws.Range('B576').Value = 'San Nicolás De Los Garza'
$ws.Range('B581').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B583').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B584').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B585').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B586').Value = 'Ixtlán De Juárez'
$ws.Range('B587').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B591').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B593').Value = 'Nejapa De Madero'
$ws.Range('B594').Value = 'Oaxaca De Juárez'
$ws.Range('B595').Value = 'Ocotlán De Morelos'
$ws.Range('B596').Value = 'Pinotepa De Don Luis'
$ws.Range('B597').Value = 'Putla Villa De Guerrero'
$ws.Range('B603').Value = 'San Dionisio Del Mar'
$ws.Range('B607').Value = 'San José Del Progreso'
$ws.Range('B619').Value = 'San Miguel Del Puerto'
$ws.Range('B626').Value = 'San Pedro El Alto'
$ws.Range('B636').Value = 'Santa Cruz Tacache De Mina'
$ws.Range('B658').Value = 'Tataltepec De Valdés'
$ws.Range('B660').Value = 'Tepelmeme Villa De Morelos'
$ws.Range('B661').Value = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
$ws.Range('B662').Value = 'Tlacolula De Matamoros'
$ws.Range('B665').Value = 'Villa De Tututepec'
$ws.Range('B666').Value = 'Villa Sola De Vega'
$ws.Range('B667').Value = 'Zimatlán De Álvarez'
$ws.Range('B681').Value = 'Cuayuca De Andrade'
$ws.Range('B687').Value = 'Huitzilan De Serdán'
$ws.Range('B690').Value = 'Izúcar De Matamoros'
$ws.Range('B691').Value = 'Los Reyes De Juárez'
$ws.Range('B695').Value = 'Palmar De Bravo'
$ws.Range('B702').Value = 'Tecali De Herrera'
$ws.Range('B706').Value = 'Tepexi De Rodríguez'
$ws.Range('B707').Value = 'Tetela De Ocampo'
$ws.Range('B710').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B715').Value = 'Totoltepec De Guerrero'
$ws.Range('B716').Value = 'Tuzamapan De Galeana'
$ws.Range('B720').Value = 'Xochitlán De Vicente Suárez'
$ws.Range('B723').Value = 'Amealco De Bonfil'
$ws.Range('B725').Value = 'Cadereyta De Montes'
$ws.Range('B730').Value = 'Jalpan De Serra'
$ws.Range('B731').Value = 'Landa De Matamoros'
$ws.Range('B733').Value = 'Pinal De Amoles'
$ws.Range('B736').Value = 'San Juan Del Río'
$ws.Range('B743').Value = 'Axtla De Terrazas'
$ws.Range('B748').Value = 'Ciudad Del Maíz'
$ws.Range('B757').Value = 'Mexquitic De Carmona'
$ws.Range('B762').Value = 'San Ciro De Acosta'
$ws.Range('B767').Value = 'Santa María Del Río'
$ws.Range('B768').Value = 'Soledad De Graciano Sánchez'
$ws.Range('B774').Value = 'Tanquián De Escobedo'
$ws.Range('B777').Value = 'Villa De Arista'
$ws.Range('B778').Value = 'Villa De Arriaga'
$ws.Range('B779').Value = 'Villa De Guadalupe'
$ws.Range('B780').Value = 'Villa De La Paz'
$ws.Range('B781').Value = 'Villa De Ramos'
$ws.Range('B782').Value = 'Villa De Reyes'
$ws.Range('B806').Value = 'Nacozari De García'
$ws.Range('B839').Value = 'Soto La Marina'
$ws.Range('B848').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B849').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B850').Value = 'Papalotla De Xicohténcatl'
$ws.Range('B860').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B864').Value = 'Boca Del Río'
$ws.Range('B866').Value = 'Castillo De Teayo'
$ws.Range('B877').Value = 'Cosamaloapan De Carpio'
$ws.Range('B882').Value = 'Hueyapan De Ocampo'
$ws.Range('B883').Value = 'Ignacio De La Llave'
$ws.Range('B891').Value = 'Martínez De La Torre'
$ws.Range('B893').Value = 'Medellín De Bravo'
$ws.Range('B897').Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range('B903').Value = 'Paso De Ovejas'
$ws.Range('B904').Value = 'Paso Del Macho'
$ws.Range('B906').Value = 'Poza Rica De Hidalgo'
$ws.Range('B911').Value = 'Sayula De Alemán'
$ws.Range('B913').Value = 'Soledad De Doblado'
$ws.Range('B917').Value = 'Tatahuicapan De Juárez'
$ws.Range('B941').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B943').Value = 'Concepción Del Oro'
$ws.Range('B953').Value = 'Jiménez Del Teul'
$ws.Range('B959').Value = 'Mezquital Del Oro'
$ws.Range('B964').Value = 'Nochistlán De Mejía'
$ws.Range('B965').Value = 'Noria De Ángeles'
$ws.Range('B975').Value = 'Teúl De González Ortega'
$ws.Range('B976').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B979').Value = 'Villa De Cos'

$ws.Rows("988:992").Delete()

"done"